$d = $word.ActiveDocument

# Locate the sentence we need to rework: "...) or using MPEI package
# installer. Source code is available on Google Code site."
$marker = ") or using MPEI package installer. Source code is available on Google Code site."
$hit = $d.Content
$found = $hit.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $s = $hit.Start

    # Keep the leading ")" exactly as-is (it is already its own correctly
    # formatted run) and clear out everything that follows it so we can
    # rebuild the remainder of the sentence piece by piece, the same way a
    # person retyping bits of the paragraph in Word would end up with
    # several small runs instead of one long one.
    $tail = $d.Range($s + 1, $hit.End)
    $tail.Text = ""

    $parts = @(
        ", using MPEI package installer or at our Google Code ",
        "site",
        ". ",
        "Source code is ",
        "also available on Google Code",
        "."
    )

    $pos = $s + 1
    foreach ($part in $parts) {
        $ins = $d.Range($pos, $pos)
        $ins.InsertBefore($part)
        $pos = $pos + $part.Length
    }
}
